$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Simple price/value updates (Price column D holds text, preserve exact formatting)
Set-TextValue "D2" "246.11"
Set-TextValue "D3" "24.16"
Set-TextValue "D4" "5.280"
Set-TextValue "D5" "0.05780"
Set-TextValue "D7" "3.144"
Set-TextValue "D8" "0.8111"
Set-TextValue "D9" "0.8608"
Set-TextValue "D10" "0.1378"
Set-TextValue "D11" "0.06988"
Set-TextValue "D12" "0.03150"
Set-TextValue "D13" "0.02924"
Set-TextValue "D14" "0.09395"
Set-TextValue "D15" "3.770"
Set-TextValue "D17" "0.04668"
Set-TextValue "D18" "0.0005999"
Set-TextValue "D19" "0.006109"
Set-TextValue "D20" "0.001235"
Set-TextValue "D21" "0.004647"

# E22 label change (drop "Worstin24h" suffix)
$ws.Range("E22").Value = "21NitroExNTX"

Set-TextValue "D24" "2.149"

# E28 label change (add "Bestin24h" suffix)
$ws.Range("E28").Value = "27UpBotsUBXTBestin24h"

Set-TextValue "D40" "0.03705"

# Rows 41-43: rotate KickToken / BKEXToken / CEJI entries
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1059"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002760"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.006299"
$ws.Range("E43").Value = "42KickTokenKICK"

Set-TextValue "D44" "0.007721"
Set-TextValue "D45" "0.00005271"

Set-TextValue "D48" "0.002405"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
